# Fruta / hortaliza, semanal
# Insert one new weekly record at row 375 (Fruta, Feria Lagunitas de Puerto
# Montt - Piña). Inserting the row shifts the existing rows 375-428 down to
# 376-429 automatically, and the sheet dimension grows from A1:T428 to
# A1:T429 - matching the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row above the current row 375, pushing everything else down.
$ws.Rows.Item(375).Insert()

# Populate the new row with the new weekly record.
$ws.Cells.Item(375, 1).Value  = 4
$ws.Cells.Item(375, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(375, 3).Value  = "Los Lagos"
$ws.Cells.Item(375, 4).Value  = 45077
$ws.Cells.Item(375, 5).Value  = 10
$ws.Cells.Item(375, 6).Value  = "Fruta"
$ws.Cells.Item(375, 7).Value  = 100108
$ws.Cells.Item(375, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(375, 9).Value  = 100108005
$ws.Cells.Item(375, 10).Value = "Pi$([char]0x00F1)a"
$ws.Cells.Item(375, 11).Value = "Caramelo"
$ws.Cells.Item(375, 12).Value = "Primera"
$ws.Cells.Item(375, 13).Value = 40
$ws.Cells.Item(375, 14).Value = 18000
$ws.Cells.Item(375, 15).Value = 19000
$ws.Cells.Item(375, 16).Value = 18500
$ws.Cells.Item(375, 17).Value = "`$/caja 12 unidades"
$ws.Cells.Item(375, 18).Value = "Ecuador"
$ws.Cells.Item(375, 19).Value = 1542
$ws.Cells.Item(375, 20).Value = 12
